$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = '{''eklabool'', ''$'', ''andamhie'', ''shimenet'', ''anda'', ''chika''}'
$ws.Range("D4").Value = '{''eklabool'', ''$'', ''andamhie'', ''shimenet'', ''chika'', ''anda'', ''naur''}'
$ws.Range("D5").Value = '{''eklabool'', ''andamhie'', ''anda'', ''chika''}'
$ws.Range("D13").Value = '{'';'', '')'', ''='', '',''}'
$ws.Range("D14").Value = '{''%'', ''!='', ''**='', ''<'', ''+='', ''//'', '';'', ''**'', '':'', ''&&'', ''||'', ''>='', '')'', ''='', ''=='', ''+'', ''step'', '']'', ''}'', ''<='', ''//='', ''%='', '','', ''-='', ''to'', ''>'', ''/='', ''/'', ''-'', ''*'', ''*=''}'
$ws.Range("D15").Value = '{''%'', ''!='', ''**='', ''<'', ''+='', ''//'', '';'', ''**'', '':'', ''&&'', ''||'', ''>='', '')'', ''='', ''=='', ''step'', ''+'', '']'', ''}'', ''<='', ''//='', ''%='', '','', ''-='', ''to'', ''>'', ''/='', ''/'', ''-'', ''*'', ''*=''}'
$ws.Range("D16").Value = '{''%'', ''!='', '']'', ''}'', ''**='', ''<='', ''//='', ''%='', '','', ''<'', ''-='', ''+='', ''//'', '';'', ''**'', ''+'', ''to'', ''>'', '':'', ''&&'', ''/='', ''||'', ''>='', '')'', ''='', ''/'', ''=='', ''*'', ''step'', ''*='', ''-''}'
$ws.Range("D17").Value = '{'';'', '']'', ''}'', '',''}'
$ws.Range("D18").Value = '{''%'', ''!='', ''<'', ''//'', '';'', ''**'', '':'', ''&&'', ''||'', ''>='', '')'', ''=='', ''+'', ''step'', '']'', ''}'', ''<='', '','', ''to'', ''>'', ''/'', ''-'', ''*''}'
$ws.Range("D19").Value = '{''%'', ''!='', '']'', ''}'', ''<='', ''step'', '','', ''<'', ''//'', '';'', ''**'', ''to'', ''>'', '':'', ''&&'', ''||'', ''>='', '')'', ''/'', ''=='', ''*'', ''+'', ''-''}'
$ws.Range("D20").Value = '{''%'', ''!='', ''<'', ''//'', '';'', ''**'', '':'', ''&&'', ''||'', ''>='', '')'', ''=='', ''step'', ''+'', '']'', ''}'', ''<='', '','', ''to'', ''>'', ''/'', ''-'', ''*''}'
$ws.Range("D21").Value = '{''%'', ''!='', ''<'', ''//'', '';'', ''**'', '':'', ''&&'', ''||'', ''>='', '')'', ''=='', ''step'', ''+'', '']'', ''}'', ''<='', '','', ''to'', ''>'', ''/'', ''-'', ''*''}'
$ws.Range("D32").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D33").Value = '{''eklabool'', ''}'', ''id'', ''ditech'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''--'', ''gogogo'', ''serve'', ''amaccana'', ''forda'', ''versa'', ''adelete'', ''chika'', ''pak''}'
$ws.Range("D37").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D39").Value = '{'';'', '']'', ''}'', ''to'', '':'', '','', '')'', ''step''}'
$ws.Range("D40").Value = '{'';'', '']'', ''}'', ''to'', '':'', '','', '')'', ''step''}'
$ws.Range("D41").Value = '{'';'', '']'', ''}'', ''to'', '':'', '','', '')'', ''step''}'
$ws.Range("D42").Value = '{'';'', '']'', ''}'', ''to'', '':'', '','', '')'', ''step''}'
$ws.Range("D43").Value = '{''%'', ''!='', ''<'', ''//'', '';'', ''**'', '':'', ''&&'', ''||'', ''>='', '')'', ''=='', ''step'', ''+'', '']'', ''}'', ''<='', '','', ''to'', ''>'', ''/'', ''-'', ''*''}'
$ws.Range("D44").Value = '{''andamhie_literal'', ''len'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''('', ''eme''}'
$ws.Range("D45").Value = '{''%'', ''!='', ''<'', ''//'', '';'', ''**'', '':'', ''&&'', ''||'', ''>='', '')'', ''=='', ''+'', ''step'', '']'', ''}'', ''<='', '','', ''to'', ''>'', ''/'', ''-'', ''*''}'
$ws.Range("D47").Value = '{''%'', ''!='', ''<'', ''//'', '';'', ''**'', '':'', ''&&'', ''||'', ''>='', '')'', ''=='', ''step'', ''+'', '']'', ''}'', ''<='', '','', ''to'', ''>'', ''/'', ''-'', ''*''}'
$ws.Range("D48").Value = '{''%'', ''!='', '']'', ''}'', ''id'', ''<='', '','', ''<'', ''//'', '';'', ''**'', ''+'', ''to'', ''>'', '':'', ''&&'', ''||'', ''>='', '')'', ''/'', ''=='', ''*'', ''step'', ''-''}'
$ws.Range("D49").Value = '{''%'', ''!='', ''<'', ''//'', '';'', ''**'', '':'', ''&&'', ''||'', ''>='', '')'', ''=='', ''+'', ''step'', '']'', ''}'', ''<='', '','', ''to'', ''>'', ''/'', ''-'', ''*''}'
$ws.Range("D50").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D51").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D52").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D54").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D55").Value = '{''andamhie_literal'', ''len'', ''!'', ''chika_literal'', ''id'', ''anda_literal'', ''--'', ''++'', ''{'', ''korik'', ''-'', ''('', ''eme''}'
$ws.Range("D57").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D60").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D61").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D62").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D64").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D67").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D69").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''ganern'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D70").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D71").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D72").Value = '{''eklabool'', ''}'', ''id'', ''ditech'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''betsung'', ''adele'', ''--'', ''andamhie'', ''gogogo'', ''serve'', ''amaccana'', ''forda'', ''versa'', ''adelete'', ''chika'', ''pak''}'
$ws.Range("D75").Value = '{'')'', ''step''}'
$ws.Range("D78").Value = '{''to'', '')'', ''step''}'
$ws.Range("D80").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D81").Value = '{''eklabool'', ''}'', ''id'', ''ditech'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''betsung'', ''adele'', ''--'', ''andamhie'', ''gogogo'', ''serve'', ''amaccana'', ''forda'', ''versa'', ''adelete'', ''chika'', ''pak''}'
$ws.Range("D83").Value = '{''eklabool'', ''}'', ''id'', ''ditech'', ''keri'', ''push'', ''++'', ''anda'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''--'', ''gogogo'', ''serve'', ''forda'', ''amaccana'', ''adelete'', ''versa'', ''chika'', ''pak''}'
$ws.Range("D84").Value = '{''}'', ''ditech''}'
$ws.Range("D86").Value = '{''amaccana'', ''}'', ''ditech'', ''betsung''}'
$ws.Range("D87").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D88").Value = '{''}'', ''ditech''}'
$ws.Range("D89").Value = '{''}'', ''ditech'', ''betsung''}'
$ws.Range("D90").Value = '{''}'', ''ditech'', ''betsung''}'
$ws.Range("D92").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D93").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
$ws.Range("D94").Value = '{''eklabool'', ''push'', ''anda'', ''--'', ''serve'', ''amaccana'', ''versa'', ''chika'', ''pak'', ''}'', ''id'', ''ditech'', ''keri'', ''++'', ''naur'', ''betsung'', ''adele'', ''andamhie'', ''gogogo'', ''forda'', ''adelete''}'
